$wb = $excel.ActiveWorkbook

# --- Create the new "Spain" sheet as a copy of "Italy" (same layout/styles),
#     placed after the last sheet (Copy with no "before" arg, only "after"
#     appends it at the end and makes it the active sheet/tab). -----------
$italy = $wb.Worksheets.Item("Italy")
$italy.Copy([System.Reflection.Missing]::Value, $italy)
$spain = $wb.Worksheets.Item($wb.Worksheets.Count)
$spain.Name = "Spain"
$spain.Activate()

# --- Make room for the three new "Attached Functionality" rows that sit
#     above the Wg/Attached Functionality-legend rows at the bottom. -------
$spain.Range("A14:A16").Insert()

# Re-apply the border-only style (row 13's) to the freshly inserted cells,
# using a format-only paste so we reuse the existing style index instead of
# fabricating a brand-new one.
$spain.Range("A13").Copy()
$spain.Range("A14:A16").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Market-specific text -------------------------------------------------
$spain.Range("B2").Value = "Spain Market"
$spain.Range("A14").Value = "PZ4DS(Dect/Fault)"
$spain.Range("A15").Value = "Three PZ4DS(Dect/Fault)"
$spain.Range("A16").Value = "Two PZ4DS(Dect/Fault)"
$spain.Range("B4").Value = "NGC-3103/T2044"

# --- Column widths re-fitted for the new (narrower) content in B:D --------
$spain.Columns.Item(2).ColumnWidth = 15.21875
$spain.Columns.Item(3).ColumnWidth = 16.6640625
$spain.Columns.Item(4).ColumnWidth = 22.33203125

# --- Row heights grow a bit once the header rows wrap over the narrower
#     columns. ---------------------------------------------------------
$spain.Range("A2:A5").EntireRow.RowHeight = 28.8

# --- Selection on the new sheet -------------------------------------------
$spain.Range("C4").Select()

Write-Host "Spain sheet added"
